$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.257.93'
$ws.Range("E2").Value = '  -0.36%  '

$ws.Range("D3").Value = '1.592.39'
$ws.Range("E3").Value = '  -0.02%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.74%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.502'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.36%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("E8").Value = '  -0.23%  '

$ws.Range("E9").Value = '  -0.58%  '

$ws.Range("E10").Value = '  -2.40%  '

$ws.Range("E11").Value = '  +0.76%  '

$ws.Range("D12").Value = '1.816.85'
$ws.Range("E12").Value = '  +0.00%  '

$ws.Range("D13").Value = '1.594.38'
$ws.Range("E13").Value = '  -0.35%  '

$ws.Range("E14").Value = '  -1.28%  '

$ws.Range("E15").Value = '  -2.50%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.89'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.08%  '

$ws.Range("D17").Value = '26.251.20'
$ws.Range("E17").Value = '  -0.41%  '

$ws.Range("E18").Value = '  -0.96%  '

$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.13'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.05%  '

$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.62%  '

$ws.Range("E21").Value = '  -0.01%  '

$ws.Range("E22").Value = '  -0.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.06'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.60%  '

$ws.Range("E24").Value = '  -4.10%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.91'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.07%  '

$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("E27").Value = '  -1.28%  '

$ws.Range("E28").Value = '  +0.22%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.11'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.58%  '

$ws.Range("E30").Value = '  -2.35%  '

$ws.Range("E31").Value = '  +0.44%  '

$ws.Range("E32").Value = '  -0.72%  '

$ws.Range("D33").Value = '1.414.98'
$ws.Range("E33").Value = '  +5.26%  '

$ws.Range("E34").Value = '  -0.16%  '

$ws.Range("E35").Value = '  -0.53%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.46'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.59%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.578'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.02%  '

$ws.Range("E38").Value = '  -1.08%  '

$ws.Range("E39").Value = '  +0.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.78'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.03%  '

$ws.Range("E41").Value = '  +0.00%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.976'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.09%  '

$ws.Range("E43").Value = '  +1.11%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.762'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.08%  '

$ws.Range("D45").Value = '1.729.07'
$ws.Range("E45").Value = '  -0.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.95'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.14%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '87.15'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.80%  '

$ws.Range("E48").Value = '  -1.07%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0510'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.97%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0953'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.50%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.01%  '
